$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 24
$ws.Range("C2").Value = 22
$ws.Range("D2").Value = 17
$ws.Range("E2").Value = 18
$ws.Range("G2").Value = 1
$ws.Range("I2").Value = 15
$ws.Range("N2").Value = 0.818
$ws.Range("O2").Value = 1.273
$ws.Range("P2").Value = 0.833
$ws.Range("Q2").Value = 2.106

# Row 3
$ws.Range("B3").Value = 23
$ws.Range("C3").Value = 21
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 13
$ws.Range("G3").Value = 1
$ws.Range("I3").Value = 11
$ws.Range("N3").Value = 0.619
$ws.Range("O3").Value = 1.143
$ws.Range("P3").Value = 0.609
$ws.Range("Q3").Value = 1.752

# Row 4
$ws.Range("B4").Value = 21
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 14
$ws.Range("G4").Value = 1
$ws.Range("I4").Value = 11
$ws.Range("N4").Value = 0.667
$ws.Range("O4").Value = 0.857
$ws.Range("P4").Value = 0.667
$ws.Range("Q4").Value = 1.524

# Row 5
$ws.Range("B5").Value = 19
$ws.Range("C5").Value = 18
$ws.Range("D5").Value = 8
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 13
$ws.Range("N5").Value = 0.667
$ws.Range("O5").Value = 1.389
$ws.Range("P5").Value = 0.632
$ws.Range("Q5").Value = 2.02

# Row 7
$ws.Range("B7").Value = 19
$ws.Range("C7").Value = 18
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 13
$ws.Range("G7").Value = 2
$ws.Range("I7").Value = 7
$ws.Range("L7").Value = 1
$ws.Range("N7").Value = 0.722
$ws.Range("O7").Value = 1.278
$ws.Range("P7").Value = 0.684
$ws.Range("Q7").Value = 1.962

# Row 9
$ws.Range("B9").Value = 18
$ws.Range("C9").Value = 18
$ws.Range("E9").Value = 9
$ws.Range("I9").Value = 4
$ws.Range("N9").Value = 0.5
$ws.Range("O9").Value = 0.5
$ws.Range("P9").Value = 0.5
$ws.Range("Q9").Value = 1

# Row 10
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 17
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 9
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 5
$ws.Range("N10").Value = 0.529
$ws.Range("O10").Value = 0.765
$ws.Range("P10").Value = 0.529
$ws.Range("Q10").Value = 1.294

# Row 12
$ws.Range("B12").Value = 19
$ws.Range("C12").Value = 17
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = 9
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 2
$ws.Range("N12").Value = 0.529
$ws.Range("O12").Value = 0.529
$ws.Range("P12").Value = 0.579
$ws.Range("Q12").Value = 1.108

# Row 13
$ws.Range("B13").Value = 19
$ws.Range("C13").Value = 18
$ws.Range("D13").Value = 9
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 1
$ws.Range("I13").Value = 6
$ws.Range("N13").Value = 0.556
$ws.Range("O13").Value = 0.611
$ws.Range("P13").Value = 0.579
$ws.Range("Q13").Value = 1.19

# Row 15
$ws.Range("B15").Value = 12
$ws.Range("C15").Value = 9
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 7
$ws.Range("L15").Value = 1
$ws.Range("N15").Value = 0.778
$ws.Range("P15").Value = 0.75
$ws.Range("Q15").Value = 1.75

# Row 18
$ws.Range("B18").Value = 239
$ws.Range("C18").Value = 224
$ws.Range("D18").Value = 101
$ws.Range("E18").Value = 141
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 7
$ws.Range("I18").Value = 101
$ws.Range("J18").Value = 11
$ws.Range("L18").Value = 4
$ws.Range("N18").Value = 0.629
$ws.Range("O18").Value = 0.902
$ws.Range("P18").Value = 0.636
$ws.Range("Q18").Value = 1.538

